$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 27454
$ws.Range("J87").Value = 27454
$ws.Range("L87").Value = 27454
$ws.Range("N87").Value = -29950
$ws.Range("H90").Value = 27454
$ws.Range("J90").Value = 27454
$ws.Range("L90").Value = 82362
$ws.Range("N90").Value = -94842
$ws.Range("H131").Value = 3697.8572
$ws.Range("I131").Value = 661.6667
$ws.Range("J131").Value = 5975
$ws.Range("K131").Value = 1985.0001
$ws.Range("L131").Value = 17925
$ws.Range("M131").Value = 3054.9999
$ws.Range("N131").Value = -28005
$ws.Range("H138").Value = 2130.158
$ws.Range("I138").Value = 1400.7273
$ws.Range("J138").Value = 3133.125
$ws.Range("K138").Value = 4202.1819
$ws.Range("L138").Value = 9399.375
$ws.Range("M138").Value = 937.8181000000004
$ws.Range("N138").Value = -19679.375
$ws.Range("H139").Value = 44925
$ws.Range("J139").Value = 44925
$ws.Range("L139").Value = 44925
$ws.Range("N139").Value = -55205
$ws.Range("H140").Value = 73250
$ws.Range("J140").Value = 73250
$ws.Range("L140").Value = 73250
$ws.Range("N140").Value = -83610

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1428.8
$ws.Range("I20").Value = 1598.5714
$ws.Range("J20").Value = 1032.6666
$ws.Range("K20").Value = 1598.5714
$ws.Range("L20").Value = 1032.6666
$ws.Range("M20").Value = -1351.5714
$ws.Range("N20").Value = -1526.6666
$ws.Range("H105").Value = 3812.9788
$ws.Range("I105").Value = 2662.7778
$ws.Range("K105").Value = 2662.7778
$ws.Range("M105").Value = -915.7777999999998
$ws.Range("H134").Value = 2552.8372
$ws.Range("I134").Value = 1404.0625
$ws.Range("J134").Value = 5894.727
$ws.Range("K134").Value = 4212.1875
$ws.Range("L134").Value = 17684.181
$ws.Range("M134").Value = -1677.1875
$ws.Range("N134").Value = -22754.181

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1037
$ws.Range("I16").Value = 1255.5
$ws.Range("J16").Value = 600
$ws.Range("K16").Value = 1255.5
$ws.Range("L16").Value = 600
$ws.Range("M16").Value = -968.5
$ws.Range("N16").Value = -1174
$ws.Range("H31").Value = 10757240
$ws.Range("I31").Value = 4860.2144
$ws.Range("J31").Value = 111112780
$ws.Range("K31").Value = 4860.2144
$ws.Range("L31").Value = 111112780
$ws.Range("M31").Value = -4565.2144
$ws.Range("N31").Value = -111113370
$ws.Range("H34").Value = 10757240
$ws.Range("I34").Value = 4860.2144
$ws.Range("J34").Value = 111112780
$ws.Range("K34").Value = 4860.2144
$ws.Range("L34").Value = 111112780
$ws.Range("M34").Value = -4658.2144
$ws.Range("N34").Value = -111113184
$ws.Range("H59").Value = 22127
$ws.Range("J59").Value = 22127
$ws.Range("L59").Value = 22127
$ws.Range("N59").Value = -24417
$ws.Range("H60").Value = 7827.25
$ws.Range("H113").Value = 1037
$ws.Range("I113").Value = 1255.5
$ws.Range("J113").Value = 600
$ws.Range("K113").Value = 1255.5
$ws.Range("L113").Value = 600
$ws.Range("M113").Value = 914.5
$ws.Range("N113").Value = -4940
$ws.Range("H140").Value = 37190
$ws.Range("J140").Value = 37190
$ws.Range("L140").Value = 37190
$ws.Range("N140").Value = -47550

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 145.3125
$ws.Range("I38").Value = 337.4
$ws.Range("J38").Value = 58
$ws.Range("K38").Value = 1012.2
$ws.Range("L38").Value = 174
$ws.Range("M38").Value = -665.1999999999999
$ws.Range("N38").Value = -868
$ws.Range("H121").Value = 1111.6666
$ws.Range("I121").Value = 577.8
$ws.Range("J121").Value = 1779
$ws.Range("K121").Value = 1733.4
$ws.Range("L121").Value = 5337
$ws.Range("M121").Value = -423.3999999999999
$ws.Range("N121").Value = -7957
$ws.Range("H122").Value = 990.2353000000001
$ws.Range("I122").Value = 859.1818
$ws.Range("K122").Value = 7732.6362
$ws.Range("M122").Value = -5282.6362

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 41996.188
$ws.Range("I70").Value = 103705
$ws.Range("J70").Value = 4970.9
$ws.Range("K70").Value = 103705
$ws.Range("L70").Value = 4970.9
$ws.Range("M70").Value = -103435
$ws.Range("N70").Value = -5510.9
$ws.Range("H73").Value = 41996.188
$ws.Range("I73").Value = 103705
$ws.Range("J73").Value = 4970.9
$ws.Range("K73").Value = 103705
$ws.Range("L73").Value = 4970.9
$ws.Range("M73").Value = -102769
$ws.Range("N73").Value = -6842.9
$ws.Range("H126").Value = 3611.3547
$ws.Range("I126").Value = 2232.9333
$ws.Range("J126").Value = 4903.625
$ws.Range("K126").Value = 6698.7999
$ws.Range("L126").Value = 14710.875
$ws.Range("M126").Value = -4228.7999
$ws.Range("N126").Value = -19650.875
$ws.Range("H138").Value = 57570.57
$ws.Range("J138").Value = 57570.57
$ws.Range("L138").Value = 57570.57
$ws.Range("N138").Value = -67850.57000000001
$ws.Range("H139").Value = 64000
$ws.Range("J139").Value = 64000
$ws.Range("L139").Value = 64000
$ws.Range("N139").Value = -74280

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4071.2856
$ws.Range("I7").Value = 4574.3125
$ws.Range("J7").Value = 3761.7307
$ws.Range("K7").Value = 4574.3125
$ws.Range("L7").Value = 3761.7307
$ws.Range("M7").Value = -4462.3125
$ws.Range("N7").Value = -3985.7307
$ws.Range("H40").Value = 3450.3057
$ws.Range("I40").Value = 4513
$ws.Range("K40").Value = 4513
$ws.Range("M40").Value = -4377
$ws.Range("H74").Value = 20197
$ws.Range("I74").Value = 20197
$ws.Range("K74").Value = 20197
$ws.Range("M74").Value = -19199
$ws.Range("H77").Value = 20197
$ws.Range("I77").Value = 20197
$ws.Range("K77").Value = 60591
$ws.Range("M77").Value = -55599
$ws.Range("H126").Value = 4071.2856
$ws.Range("I126").Value = 4574.3125
$ws.Range("J126").Value = 3761.7307
$ws.Range("K126").Value = 13722.9375
$ws.Range("L126").Value = 11285.1921
$ws.Range("M126").Value = -11252.9375
$ws.Range("N126").Value = -16225.1921
$ws.Range("H139").Value = 41995.4
$ws.Range("J139").Value = 42144.89
$ws.Range("L139").Value = 42144.89
$ws.Range("N139").Value = -52424.89
